$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeInfo")
$ws.Range("D9").Value = "suvarna"
$ws.Range("E7").Select()
